$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 2.5
$ws.Range("A11").Value = 2.6
$ws.Range("A12").Value = 2.7

$ws.Range("A13").Select()
